# Generate Report for Handoff
#
# Updates the localization-status report to reflect that the zh-cn and
# de-de handoff packages have moved from "In Translation" to
# "Ready for handoff", refreshes the corresponding handoff timestamps,
# and widens the Status column(s) to fit the new, longer status text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---- Overview sheet -------------------------------------------------
# E2/F2 hold the per-language status ("In Translation" -> "Ready for handoff")
# G2 holds the "Latest HO Xliff Generate Date" timestamp.
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-24 17:01:23"

# ---- zh-cn detail sheet ----------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-24 17:01:18"

# ---- de-de detail sheet ----------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-24 17:01:23"

# ---- Column widths -----------------------------------------------------
# The longer "Ready for handoff" text made the Status columns resize
# (equivalent to the columns being auto-fit to the new content).
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332
$wsZhCn.Columns.Item(3).ColumnWidth     = 16.333333333333332
$wsDeDe.Columns.Item(3).ColumnWidth     = 16.333333333333332
